# Insert a new data row at row 854 (pushing the existing rows 854-931 down
# to 855-932) and populate it with the new record described by the diff.
# The workbook has a single sheet; row 1 is the header, data starts at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 854. Excel shifts row 854..931
# down to 855..932 and the new blank row inherits formatting (incl. the
# date-style on column D) from the row that used to be there.
$ws.Rows.Item(854).Insert()

# Fill in the new row 854 with its values. Columns A, B, C, E, F, G, H, I,
# Q, R are identical to the row that used to occupy this slot (now at 855),
# so Excel's Insert already duplicated that formatting/content context and
# we only need to set the literal values explicitly for every column of
# the new record.
$ws.Cells.Item(854, 1).Value = 6
$ws.Cells.Item(854, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(854, 3).Value = "Metropolitana"
$ws.Cells.Item(854, 4).Value = 44769
$ws.Cells.Item(854, 5).Value = 13
$ws.Cells.Item(854, 6).Value = 100112031
$ws.Cells.Item(854, 7).Value = "Poroto verde"
$ws.Cells.Item(854, 8).Value = "Magnum"
$ws.Cells.Item(854, 9).Value = "Primera"
$ws.Cells.Item(854, 10).Value = 220
$ws.Cells.Item(854, 11).Value = 28000
$ws.Cells.Item(854, 12).Value = 30000
$ws.Cells.Item(854, 13).Value = 28909
$ws.Cells.Item(854, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(854, 15).Value = "Perú"
$ws.Cells.Item(854, 16).Value = 1156
$ws.Cells.Item(854, 17).Value = 25
$ws.Cells.Item(854, 18).Value = "Hortaliza"
